# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Macroferia Regional de Talca - Frutilla"
# at the top of the data block that starts at row 544, pushing the existing
# rows 544:570 down to 547:573 (dimension grows from A1:T570 to A1:T573).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 544; this shifts the existing
# rows 544-570 down to 547-573 and extends the used range automatically.
$ws.Rows("544:546").Insert()

# Columns A, B, C, E-K and T are constant across this whole data block, so
# copy them down into the three freshly inserted rows.
$commonCols = @("A", "B", "C", "E", "F", "G", "H", "I", "J", "K", "T")
foreach ($col in $commonCols) {
    $src = $ws.Range($col + "547").Value2
    $ws.Range($col + "544").Value = $src
    $ws.Range($col + "545").Value = $src
    $ws.Range($col + "546").Value = $src
}

# Row 544: Especial
$ws.Range("D544").Value = 44706
$ws.Range("L544").Value = "Especial"
$ws.Range("M544").Value = 50
$ws.Range("N544").Value = 13000
$ws.Range("O544").Value = 13000
$ws.Range("P544").Value = 13000
$ws.Range("Q544").Value = "$/bandeja 7 kilos"
$ws.Range("R544").Value = "Provincia de Melipilla"
$ws.Range("S544").Value = 1857

# Row 545: Primera
$ws.Range("D545").Value = 44706
$ws.Range("L545").Value = "Primera"
$ws.Range("M545").Value = 30
$ws.Range("N545").Value = 11000
$ws.Range("O545").Value = 11000
$ws.Range("P545").Value = 11000
$ws.Range("Q545").Value = "$/bandeja 7 kilos"
$ws.Range("R545").Value = "Provincia de Melipilla"
$ws.Range("S545").Value = 1571

# Row 546: Segunda
$ws.Range("D546").Value = 44706
$ws.Range("L546").Value = "Segunda"
$ws.Range("M546").Value = 20
$ws.Range("N546").Value = 5000
$ws.Range("O546").Value = 5000
$ws.Range("P546").Value = 5000
$ws.Range("Q546").Value = "$/bandeja 7 kilos"
$ws.Range("R546").Value = "Provincia de Melipilla"
$ws.Range("S546").Value = 714

# Apply the same date number format used by the rest of column D to the
# newly created date cells.
$dateFmt = $ws.Range("D547").NumberFormat
$ws.Range("D544:D546").NumberFormat = $dateFmt
